$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Range("L5").Value = 5273.85
$ws1.Range("L15").Value = 3873.05
$ws1.Range("M34").Value = 2673.62
$ws1.Range("P34").Value = 1651.9
$ws1.Range("M48").Value = 1676.27
$ws1.Range("M49").Value = 5765.94
$ws1.Range("P60").Value = "5 de 58"

$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Range("F5").Value = 21118.81
$ws2.Range("F15").Value = 8398.709999999999
$ws2.Range("F34").Value = 10282.96
$ws2.Range("F48").Value = 3626.8
$ws2.Range("F49").Value = 6664.04
$ws2.Range("F60").Value = 84583.83

$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$ws3.Range("D8").Value = 3081.08
$ws3.Range("E8").Value = -2632.27837082797
$ws3.Range("F8").Value = 6.865126594313213

$ws3.Range("D11").Value = 20382.05
$ws3.Range("E11").Value = -808.989750750301
$ws3.Range("F11").Value = 1.041331796890643

$ws3.Range("D12").Value = 43953.03
$ws3.Range("E12").Value = 4671.029999999999
$ws3.Range("F12").Value = 0.9039358292993223

$ws3.Range("D14").Value = 90187.98999999999
$ws3.Range("E14").Value = 9710.002841887861
$ws3.Range("F14").Value = 0.9028008214613858
